$d = $word.ActiveDocument

# The document's final paragraph is an (otherwise empty) paragraph that
# only carries the _GoBack bookmark. Insert the three new minutes
# paragraphs immediately before it (as a single InsertBefore call so
# they land in document order), leaving the bookmark paragraph last.
$last = $d.Paragraphs.Last
$r = $last.Range

$newText = "Installed ICLOCS, getting familiar with it.`r" + `
    "Problem formulation for single agent QUAV: written down, now must implement using ICLOCS.`r" + `
    "Next: attempt to implement open-loop control, then closed-loop control with Simulink, then CLC with ROS.`r"

$r.InsertBefore($newText)
